$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D5, D7, D8 to the text value "Free" (these cells become shared-string text cells)
$ws.Range("D5").Value = "Free"
$ws.Range("D7").Value = "Free"
$ws.Range("D8").Value = "Free"

# Update the visible selection to match the saved view state
$ws.Range("D2").Select()
